# Edgar_scaling_mapping.xlsx update: switch to EDGAR v4.3.2 for all species (incl. CH4)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "map": update / add Edgar_Sector_Name values (column B) for several
# inv_sector rows, matching the new EDGAR sector naming.
# ---------------------------------------------------------------------------
$mapWs = $wb.Worksheets.Item("map")

$mapWs.Range("B29").Value = "Fugitive emissions from oil and gas"
$mapWs.Range("B38").Value = "Production of metals"
$mapWs.Range("B39").Value = "Production of pulp/paper/food/drink"
$mapWs.Range("B52").Value = "Agricultural waste burning"
$mapWs.Range("B60").Value = "Fossil fuel fires"

# ---------------------------------------------------------------------------
# Sheet "year": EDGAR inventory now extends through 2012 (was 2009) for the
# "end_scaling_year" column (G) on all rows that previously read 2009.
# ---------------------------------------------------------------------------
$yearWs = $wb.Worksheets.Item("year")

$yearWs.Range("G2").Value = 2012
$yearWs.Range("G3").Value = 2012
$yearWs.Range("G8:G31").Value = 2012

# ---------------------------------------------------------------------------
# Restore the "year" sheet as the active/selected tab (was "map" before).
# ---------------------------------------------------------------------------
$yearWs.Activate()
